$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Cells.Item(4, 6).Value = 470
$wsExhibit.Cells.Item(6, 6).Value = 225
$wsExhibit.Cells.Item(7, 6).Value = 208
$wsExhibit.Cells.Item(8, 6).Value = 239
$wsExhibit.Cells.Item(9, 6).Value = 2824
$wsExhibit.Cells.Item(10, 6).Value = 56
$wsExhibit.Cells.Item(11, 6).Value = 118
$wsExhibit.Cells.Item(12, 6).Value = 2163
$wsExhibit.Cells.Item(13, 6).Value = 253
$wsExhibit.Cells.Item(16, 6).Value = 77
$wsExhibit.Cells.Item(17, 6).Value = 2512
$wsExhibit.Cells.Item(19, 6).Value = 1247
$wsExhibit.Cells.Item(20, 6).Value = 4531
$wsExhibit.Cells.Item(22, 6).Value = 4765
$wsExhibit.Cells.Item(23, 6).Value = 1300
$wsExhibit.Cells.Item(24, 6).Value = 2770
$wsExhibit.Cells.Item(25, 6).Value = 3177
$wsExhibit.Cells.Item(27, 6).Value = 1481
$wsExhibit.Cells.Item(28, 6).Value = 232
$wsExhibit.Cells.Item(31, 6).Value = 240
$wsExhibit.Cells.Item(32, 6).Value = 850
$wsExhibit.Cells.Item(33, 6).Value = 1472
$wsExhibit.Cells.Item(34, 6).Value = 109
$wsExhibit.Cells.Item(35, 6).Value = 232
$wsExhibit.Cells.Item(36, 6).Value = 583
$wsExhibit.Cells.Item(38, 6).Value = 278
$wsExhibit.Cells.Item(39, 6).Value = 355

$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Cells.Item(3, 6).Value = 90

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Cells.Item(4, 6).Value = 470
$wsAll.Cells.Item(5, 6).Value = 90
$wsAll.Cells.Item(8, 6).Value = 225
$wsAll.Cells.Item(9, 6).Value = 208
$wsAll.Cells.Item(11, 6).Value = 239
$wsAll.Cells.Item(12, 6).Value = 2824
$wsAll.Cells.Item(13, 6).Value = 56
$wsAll.Cells.Item(14, 6).Value = 118
$wsAll.Cells.Item(15, 6).Value = 2163
$wsAll.Cells.Item(16, 6).Value = 253
$wsAll.Cells.Item(19, 6).Value = 77
$wsAll.Cells.Item(21, 6).Value = 2512
$wsAll.Cells.Item(22, 6).Value = 1247
$wsAll.Cells.Item(26, 6).Value = 4531
$wsAll.Cells.Item(28, 6).Value = 4765
$wsAll.Cells.Item(29, 6).Value = 1300
$wsAll.Cells.Item(30, 6).Value = 2770
$wsAll.Cells.Item(31, 6).Value = 3177
$wsAll.Cells.Item(35, 6).Value = 1481
$wsAll.Cells.Item(37, 6).Value = 232
$wsAll.Cells.Item(40, 6).Value = 240
$wsAll.Cells.Item(41, 6).Value = 850
$wsAll.Cells.Item(43, 6).Value = 1472
$wsAll.Cells.Item(44, 6).Value = 109
$wsAll.Cells.Item(45, 6).Value = 232
$wsAll.Cells.Item(46, 6).Value = 583
$wsAll.Cells.Item(48, 6).Value = 278
$wsAll.Cells.Item(49, 6).Value = 355
